$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the B and C columns for the header row + existing data rows (1..19) ---
# Column B currently holds the Date (number format "YYYY-MM-DD HH:MM:SS" for rows
# 2-18, "YYYY-MM-DD" for the last row 19), column C holds the price (General format).
# After the edit: column B holds the price (General) and column C holds the Date
# ("YYYY-MM-DD HH:MM:SS" for every row, since the special last-row date format now
# belongs to the new last row, 20).
for ($r = 1; $r -le 19; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2

    # write swapped values
    $cCell.Value2 = $bVal
    $bCell.Value2 = $cVal
}

# Column B (now prices) goes back to plain/General formatting for every data row.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).ClearFormats()
}

# Column C (now dates) gets the regular date/time format for every data row,
# including row 19 which used to carry the special "last row" format.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Append the new row 20, which becomes the new "last row" ---
$ws.Cells.Item(20, 1).Value2 = 781.86

$ws.Cells.Item(20, 2).Value2 = 679.38

$ws.Cells.Item(20, 3).Value2 = 45754
$ws.Cells.Item(20, 3).NumberFormat = "YYYY-MM-DD"
